$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# GALs quantity changed from 3 to 6
$ws.Range("C6").Value = 6

# Row 9: was "Goat" / "Because smoke demons" -> now "SRAM" entry
$ws.Range("B9").Value = "SRAM"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = "Store measurment data"
$ws.Range("G9").Value = "have them"

# Update selection to match the diff (cosmetic, matches author's last cursor position)
$ws.Range("D37").Select()
